# Add 2022-Q1 data:
#  - the existing "总计" sheet becomes the new "2022-Q1" sheet (with per-fund
#    holding detail, same layout as the other quarterly sheets)
#  - a brand-new "总计" sheet is appended at the end with the historical
#    summary table, now including the 2022-Q1 row at the top

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# style donor: the "2021-Q4" sheet uses the exact same layout/styling that
# the new per-fund detail sheet needs (bold/boxed header row + boxed index column)
$donor = $wb.Worksheets.Item("2021-Q4")

$donor.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$donor.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# index column (A2:A6) -- plain numbers 0..4
for ($i = 0; $i -lt 5; $i++) {
    $q1.Cells.Item(2 + $i, 1).Value = $i
}

# fund holding detail rows. Columns B-G are stored as text (even though some
# look numeric) to match the source data / avoid losing things such as
# leading zeros in fund codes; column H is a genuine number.
$rows = @(
    @("001304", "建信鑫安回报灵活配置混合", "2.13", "66.83", "6.72", "0.1431", 2),
    @("515300", "嘉实沪深300红利低波动ETF", "0.87", "99.20", "3.06", "0.0266", 8),
    @("012977", "瑞达鑫红量化6个月持有混合型证券投资基金A", "1.04", "94.56", "1.62", "0.0168", 3),
    @("510290", "南方上证380ETF", "1.75", "99.12", "0.93", "0.0163", 4),
    @("012978", "瑞达鑫红量化6个月持有混合型证券投资基金C", "0.17", "94.56", "1.62", "0.0028", 3)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowIndex = 2 + $r

    $codeCell = $q1.Cells.Item($rowIndex, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $q1.Cells.Item($rowIndex, 3).Value = $row[1]

    $sizeCell = $q1.Cells.Item($rowIndex, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[2]

    $posCell = $q1.Cells.Item($rowIndex, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[3]

    $ratioCell = $q1.Cells.Item($rowIndex, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $row[4]

    $valueCell = $q1.Cells.Item($rowIndex, 7)
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $row[5]

    $q1.Cells.Item($rowIndex, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the updated summary table
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$donor.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$donor.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$summary = @(
    @("2022-Q1", 5, 0.21),
    @("2021-Q4", 5, 0.26),
    @("2021-Q3", 1, 0.02),
    @("2021-Q2", 4, 0.05),
    @("2021-Q1", 5, 0.17),
    @("2020-Q4", 3, 0.02)
)

for ($r = 0; $r -lt $summary.Length; $r++) {
    $row = $summary[$r]
    $rowIndex = 2 + $r
    $total.Cells.Item($rowIndex, 1).Value = $r
    $total.Cells.Item($rowIndex, 2).Value = $row[0]
    $total.Cells.Item($rowIndex, 3).Value = $row[1]
    $total.Cells.Item($rowIndex, 4).Value = $row[2]
}

$total.Range("A1").Select()
